$wb = $excel.ActiveWorkbook

# Reference a cell that already carries the bold+bordered "header" style (style index 1)
# used throughout this workbook for header rows and the leading index column.
$ws1 = $wb.Worksheets.Item(1)
$headerStyleSource = $ws1.Cells.Item(1, 2)

function Set-HeaderStyle($range) {
    $headerStyleSource.Copy()
    $range.PasteSpecial(-4122)
}

# ===== sheet "保險" (insurance) =====
$ws6 = $wb.Worksheets.Item("保險")
$ws6.Cells.Clear()
$ws6.Range("B1").Value = "company"
$ws6.Range("C1").Value = "name"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "property_category"
$ws6.Range("F1").Value = "category"
$ws6.Range("G1").Value = "date"
$ws6.Range("H1").Value = "legislator_name"
$ws6.Range("I1").Value = "legislator_id"
$ws6.Range("J1").Value = "source_file"
$ws6.Range("K1").Value = "index"
$ws6.Range("A2").Value = 89
$ws6.Range("B2").Value = "富邦人壽"
$ws6.Range("C2").Value = "20LPL安泰分紅終身壽險"
$ws6.Range("D2").Value = "蘇震清"
$ws6.Range("E2").Value = "insurance"
$ws6.Range("F2").Value = "normal"
$ws6.Range("G2").Value = "2011-11-17"
$ws6.Range("H2").Value = "蘇震清"
$ws6.Range("I2").Value = 1718
$ws6.Range("J2").Value = "tmp98701"
$ws6.Range("K2").Value = 89
$ws6.Range("A3").Value = 90
$ws6.Range("B3").Value = "富邦人壽"
$ws6.Range("C3").Value = "20LPL安泰分红終身壽險"
$ws6.Range("D3").Value = "廖靖汝"
$ws6.Range("E3").Value = "insurance"
$ws6.Range("F3").Value = "normal"
$ws6.Range("G3").Value = "2011-11-17"
$ws6.Range("H3").Value = "蘇震清"
$ws6.Range("I3").Value = 1718
$ws6.Range("J3").Value = "tmp98701"
$ws6.Range("K3").Value = 90
$ws6.Range("A4").Value = 91
$ws6.Range("B4").Value = "富邦人壽"
$ws6.Range("C4").Value = "20LPL安泰分紅終身壽險"
$ws6.Range("D4").Value = "蘇〇淳"
$ws6.Range("E4").Value = "insurance"
$ws6.Range("F4").Value = "normal"
$ws6.Range("G4").Value = "2011-11-17"
$ws6.Range("H4").Value = "蘇震清"
$ws6.Range("I4").Value = 1718
$ws6.Range("J4").Value = "tmp98701"
$ws6.Range("K4").Value = 91
$ws6.Range("A5").Value = 92
$ws6.Range("B5").Value = "富邦人壽"
$ws6.Range("C5").Value = "20LPL安泰分紅終身壽險"
$ws6.Range("D5").Value = "蘇〇婕"
$ws6.Range("E5").Value = "insurance"
$ws6.Range("F5").Value = "normal"
$ws6.Range("G5").Value = "2011-11-17"
$ws6.Range("H5").Value = "蘇震清"
$ws6.Range("I5").Value = 1718
$ws6.Range("J5").Value = "tmp98701"
$ws6.Range("K5").Value = 92
$ws6.Range("A6").Value = 93
$ws6.Range("B6").Value = "富邦人壽"
$ws6.Range("C6").Value = "安泰喬壽還本終身壽險"
$ws6.Range("D6").Value = "蘇〇淳"
$ws6.Range("E6").Value = "insurance"
$ws6.Range("F6").Value = "normal"
$ws6.Range("G6").Value = "2011-11-17"
$ws6.Range("H6").Value = "蘇震清"
$ws6.Range("I6").Value = 1718
$ws6.Range("J6").Value = "tmp98701"
$ws6.Range("K6").Value = 93
$ws6.Range("A7").Value = 94
$ws6.Range("B7").Value = "富邦人壽"
$ws6.Range("C7").Value = "安泰喬壽還本終身壽險"
$ws6.Range("D7").Value = "蘇〇婕"
$ws6.Range("E7").Value = "insurance"
$ws6.Range("F7").Value = "normal"
$ws6.Range("G7").Value = "2011-11-17"
$ws6.Range("H7").Value = "蘇震清"
$ws6.Range("I7").Value = 1718
$ws6.Range("J7").Value = "tmp98701"
$ws6.Range("K7").Value = 94
Set-HeaderStyle $ws6.Range("B1:K1")
Set-HeaderStyle $ws6.Range("A2:A7")

# ===== sheet "債務" (debt) =====
$ws7 = $wb.Worksheets.Item("債務")
$ws7.Cells.Clear()
$ws7.Range("B1").Value = "species"
$ws7.Range("C1").Value = "debtor"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"
$ws7.Range("A2").Value = 104
$ws7.Range("B2").Value = "綜合存款存摺(透支)."
$ws7.Range("C2").Value = "蘇震清"
$ws7.Range("D2").Value = "台灣銀行屏東分行屏東縣屏東市中山路"
$ws7.Range("E2").Value = 5339618
$ws7.Range("F2").Value = "100年03月28日"
$ws7.Range("G2").Value = "貸款一中期擔保放款"
$ws7.Range("H2").Value = "debt"
$ws7.Range("I2").Value = "normal"
$ws7.Range("J2").Value = "2011-11-17"
$ws7.Range("K2").Value = "蘇震清"
$ws7.Range("L2").Value = 1718
$ws7.Range("M2").Value = "tmp98701"
$ws7.Range("N2").Value = 104
$ws7.Range("A3").Value = 105
$ws7.Range("B3").Value = "長期擔保放款"
$ws7.Range("C3").Value = "廖靖汝"
$ws7.Range("D3").Value = "台灣土地銀行屏東縣屏東市逢甲路"
$ws7.Range("E3").Value = 663000
$ws7.Range("F3").Value = "87年11月03日"
$ws7.Range("G3").Value = "房貸"
$ws7.Range("H3").Value = "debt"
$ws7.Range("I3").Value = "normal"
$ws7.Range("J3").Value = "2011-11-17"
$ws7.Range("K3").Value = "蘇震清"
$ws7.Range("L3").Value = 1718
$ws7.Range("M3").Value = "tmp98701"
$ws7.Range("N3").Value = 105
$ws7.Range("A4").Value = 106
$ws7.Range("B4").Value = "綜合存款存摺(透支）"
$ws7.Range("C4").Value = "蘇震清"
$ws7.Range("D4").Value = "台灣銀行屏東分行屏東縣屏東市中山路"
$ws7.Range("E4").Value = 371101
$ws7.Range("F4").Value = "100年03月28日"
$ws7.Range("G4").Value = "貸款一治家成長貸款"
$ws7.Range("H4").Value = "debt"
$ws7.Range("I4").Value = "normal"
$ws7.Range("J4").Value = "2011-11-17"
$ws7.Range("K4").Value = "蘇震清"
$ws7.Range("L4").Value = 1718
$ws7.Range("M4").Value = "tmp98701"
$ws7.Range("N4").Value = 106
Set-HeaderStyle $ws7.Range("B1:N1")
Set-HeaderStyle $ws7.Range("A2:A4")

# ===== new sheet "事業投資" (business investment) appended at the end =====
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "事業投資"
$ws8.Range("B1").Value = "owner"
$ws8.Range("C1").Value = "company"
$ws8.Range("D1").Value = "address"
$ws8.Range("E1").Value = "total"
$ws8.Range("F1").Value = "register_date"
$ws8.Range("G1").Value = "register_reason"
$ws8.Range("H1").Value = "property_category"
$ws8.Range("I1").Value = "category"
$ws8.Range("J1").Value = "date"
$ws8.Range("K1").Value = "legislator_name"
$ws8.Range("L1").Value = "legislator_id"
$ws8.Range("M1").Value = "source_file"
$ws8.Range("N1").Value = "index"
$ws8.Range("A2").Value = 111
$ws8.Range("B2").Value = "寧靖汝"
$ws8.Range("C2").Value = "南島休閒育樂股份有限公司"
$ws8.Range("D2").Value = "高雄市精富路148號"
$ws8.Range("E2").Value = 1000000
$ws8.Range("F2").Value = "95年08月23日"
$ws8.Range("G2").Value = "投資"
$ws8.Range("H2").Value = "investment"
$ws8.Range("I2").Value = "normal"
$ws8.Range("J2").Value = "2011-11-17"
$ws8.Range("K2").Value = "蘇震清"
$ws8.Range("L2").Value = 1718
$ws8.Range("M2").Value = "tmp98701"
$ws8.Range("N2").Value = 111
Set-HeaderStyle $ws8.Range("B1:N1")
Set-HeaderStyle $ws8.Range("A2:A2")

